$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(539, 1).Value = 46028
$ws.Cells.Item(539, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(539, 2).Value = "Indiana Pacers"
$ws.Cells.Item(539, 3).Value = "Cleveland Cavaliers"
$ws.Cells.Item(539, 4).Value = 6.5
$ws.Cells.Item(539, 5).Value = 116
$ws.Cells.Item(539, 6).Value = 120
$ws.Cells.Item(539, 7).Value = 2.5

$ws.Cells.Item(540, 1).Value = 46028
$ws.Cells.Item(540, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(540, 2).Value = "Washington Wizards"
$ws.Cells.Item(540, 3).Value = "Orlando Magic"
$ws.Cells.Item(540, 4).Value = 7.5
$ws.Cells.Item(540, 5).Value = 120
$ws.Cells.Item(540, 6).Value = 112
$ws.Cells.Item(540, 7).Value = 15.5

$ws.Cells.Item(541, 1).Value = 46028
$ws.Cells.Item(541, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(541, 2).Value = "New Orleans Pelicans"
$ws.Cells.Item(541, 3).Value = "Los Angeles Lakers"
$ws.Cells.Item(541, 4).Value = 5.5
$ws.Cells.Item(541, 5).Value = 103
$ws.Cells.Item(541, 6).Value = 111
$ws.Cells.Item(541, 7).Value = -2.5

$ws.Cells.Item(542, 1).Value = 46028
$ws.Cells.Item(542, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(542, 2).Value = "Memphis Grizzlies"
$ws.Cells.Item(542, 3).Value = "San Antonio Spurs"
$ws.Cells.Item(542, 4).Value = 5.5
$ws.Cells.Item(542, 5).Value = 106
$ws.Cells.Item(542, 6).Value = 105
$ws.Cells.Item(542, 7).Value = 6.5

$ws.Cells.Item(543, 1).Value = 46028
$ws.Cells.Item(543, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(543, 2).Value = "Minnesota Timberwolves"
$ws.Cells.Item(543, 3).Value = "Miami Heat"
$ws.Cells.Item(543, 4).Value = -5.5
$ws.Cells.Item(543, 5).Value = 122
$ws.Cells.Item(543, 6).Value = 94
$ws.Cells.Item(543, 7).Value = 22.5

$ws.Cells.Item(544, 1).Value = 46028
$ws.Cells.Item(544, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(544, 2).Value = "Sacramento Kings"
$ws.Cells.Item(544, 3).Value = "Dallas Mavericks"
$ws.Cells.Item(544, 4).Value = 4.5
$ws.Cells.Item(544, 5).Value = 98
$ws.Cells.Item(544, 6).Value = 100
$ws.Cells.Item(544, 7).Value = 2.5

$ws.Cells.Item(545, 1).Value = 46029
$ws.Cells.Item(545, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(545, 2).Value = "Detroit Pistons"
$ws.Cells.Item(545, 3).Value = "Chicago Bulls"
$ws.Cells.Item(545, 4).Value = -7.5
$ws.Cells.Item(545, 5).Value = 108
$ws.Cells.Item(545, 6).Value = 93
$ws.Cells.Item(545, 7).Value = 7.5

$ws.Cells.Item(546, 1).Value = 46029
$ws.Cells.Item(546, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(546, 2).Value = "Philadelphia 76ers"
$ws.Cells.Item(546, 3).Value = "Washington Wizards"
$ws.Cells.Item(546, 4).Value = -16.5
$ws.Cells.Item(546, 5).Value = 131
$ws.Cells.Item(546, 6).Value = 110
$ws.Cells.Item(546, 7).Value = 4.5

$ws.Cells.Item(547, 1).Value = 46029
$ws.Cells.Item(547, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(547, 2).Value = "Charlotte Hornets"
$ws.Cells.Item(547, 3).Value = "Toronto Raptors"
$ws.Cells.Item(547, 4).Value = 2.5
$ws.Cells.Item(547, 5).Value = 96
$ws.Cells.Item(547, 6).Value = 97
$ws.Cells.Item(547, 7).Value = 1.5

$ws.Cells.Item(548, 1).Value = 46029
$ws.Cells.Item(548, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(548, 2).Value = "Boston Celtics"
$ws.Cells.Item(548, 3).Value = "Denver Nuggets"
$ws.Cells.Item(548, 4).Value = -10.5
$ws.Cells.Item(548, 5).Value = 110
$ws.Cells.Item(548, 6).Value = 114
$ws.Cells.Item(548, 7).Value = -14.5

$ws.Cells.Item(549, 1).Value = 46029
$ws.Cells.Item(549, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(549, 2).Value = "Atlanta Hawks"
$ws.Cells.Item(549, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(549, 4).Value = -11.5
$ws.Cells.Item(549, 5).Value = 117
$ws.Cells.Item(549, 6).Value = 100
$ws.Cells.Item(549, 7).Value = 5.5

$ws.Cells.Item(550, 1).Value = 46029
$ws.Cells.Item(550, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(550, 2).Value = "Brooklyn Nets"
$ws.Cells.Item(550, 3).Value = "Orlando Magic"
$ws.Cells.Item(550, 4).Value = 1.5
$ws.Cells.Item(550, 5).Value = 103
$ws.Cells.Item(550, 6).Value = 104
$ws.Cells.Item(550, 7).Value = 0.5

$ws.Cells.Item(551, 1).Value = 46029
$ws.Cells.Item(551, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(551, 2).Value = "New York Knicks"
$ws.Cells.Item(551, 3).Value = "Los Angeles Clippers"
$ws.Cells.Item(551, 4).Value = -5.5
$ws.Cells.Item(551, 5).Value = 123
$ws.Cells.Item(551, 6).Value = 111
$ws.Cells.Item(551, 7).Value = 6.5

$ws.Cells.Item(552, 1).Value = 46029
$ws.Cells.Item(552, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(552, 2).Value = "Memphis Grizzlies"
$ws.Cells.Item(552, 3).Value = "Phoenix Suns"
$ws.Cells.Item(552, 4).Value = 5.5
$ws.Cells.Item(552, 5).Value = 98
$ws.Cells.Item(552, 6).Value = 117
$ws.Cells.Item(552, 7).Value = -13.5

$ws.Cells.Item(553, 1).Value = 46029
$ws.Cells.Item(553, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(553, 2).Value = "Oklahoma City Thunder"
$ws.Cells.Item(553, 3).Value = "Utah Jazz"
$ws.Cells.Item(553, 4).Value = -19.5
$ws.Cells.Item(553, 5).Value = 129
$ws.Cells.Item(553, 6).Value = 125
$ws.Cells.Item(553, 7).Value = -15.5

$ws.Cells.Item(554, 1).Value = 46029
$ws.Cells.Item(554, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(554, 2).Value = "San Antonio Spurs"
$ws.Cells.Item(554, 3).Value = "Los Angeles Lakers"
$ws.Cells.Item(554, 4).Value = -8.5
$ws.Cells.Item(554, 5).Value = 107
$ws.Cells.Item(554, 6).Value = 91
$ws.Cells.Item(554, 7).Value = 7.5

$ws.Cells.Item(555, 1).Value = 46029
$ws.Cells.Item(555, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(555, 2).Value = "Golden State Warriors"
$ws.Cells.Item(555, 3).Value = "Milwaukee Bucks"
$ws.Cells.Item(555, 4).Value = -6.5
$ws.Cells.Item(555, 5).Value = 120
$ws.Cells.Item(555, 6).Value = 113
$ws.Cells.Item(555, 7).Value = 0.5

$ws.Cells.Item(556, 1).Value = 46029
$ws.Cells.Item(556, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(556, 2).Value = "Portland Trail Blazers"
$ws.Cells.Item(556, 3).Value = "Houston Rockets"
$ws.Cells.Item(556, 4).Value = 7.5
$ws.Cells.Item(556, 5).Value = 103
$ws.Cells.Item(556, 6).Value = 102
$ws.Cells.Item(556, 7).Value = 8.5
